$wb = $excel.ActiveWorkbook
$wsContacts = $wb.Worksheets.Item("Contacts")
$ws = $wb.Worksheets.Item("TestCase")

# --- Step 1: Wording corrections (click -> Click, type -> Enter, etc.) ---
$ws.Range("E7").Value = 'Click "Welcome Paul" link at the right top corener of the page'
$ws.Range("E8").Value = 'Click "Logout" button'
$ws.Range("E10").Value = 'Click "Forgot your password?" link'
$ws.Range("E11").Value = 'Input valid username into "OrangeHRM Username" input box'
$ws.Range("E12").Value = 'Click "Reset Password" button'
$ws.Range("E28").Value = 'Click "Forgot your password?" link'
$ws.Range("E29").Value = 'Click "Reset Password" button'
$ws.Range("E31").Value = 'Click "Forgot your password?" link'
$ws.Range("E32").Value = 'Input special characters in the "OrangeHRM Username" input box '
$ws.Range("E33").Value = 'Click "Reset Password" button'
$ws.Range("E37").Value = 'Navigate to MyInfo Page'
$ws.Range("E39").Value = 'Enter "Mily" in the firstname inputbox'
$ws.Range("F39").Value = 'User shuold be able to enter "Mily" in the inputbox'
$ws.Range("E40").Value = 'Enter "Hm" in the lastname inputbox'
$ws.Range("F40").Value = 'User should be able to enter "Hm" in the inputbox'
$ws.Range("E41").Value = 'Enter "007" in the "Employee ID" inputbox'
$ws.Range("F41").Value = 'User should be able to enter "007" in the "Employee ID" inputbox'
$ws.Range("E43").Value = 'Click "Female" radioButton in "Gender " '
$ws.Range("E52").Value = 'Enter file path'
$ws.Range("F52").Value = 'User should be able to enter file path'
$ws.Range("E57").Value = 'Click "Admin" link'
$ws.Range("F57").Value = 'User should see following links: "User Management","Job","Organization","Qualifications","Nationalities","Configuration";
'
$ws.Range("E58").Value = 'Navigate to admin Page'
$ws.Range("E59").Value = 'Click "Add" button'
$ws.Range("E60").Value = 'Click "User Role" label'
$ws.Range("E62").Value = 'Enter "Mily Hm" in hte "Employee Name" textbox'
$ws.Range("E63").Value = 'Enter "Mily." in the username textbox'
$ws.Range("E66").Value = 'Enter "Mily1234" in the password textbox'
$ws.Range("E68").Value = 'Click save button'
$ws.Range("E69").Value = 'Verify "Mily" display in the list '
$ws.Range("D70").Value = 'User should be able to seach a system user from admin page'
$ws.Range("E71").Value = 'Enter "Milyyy"in the "username" textField'
$ws.Range("F71").Value = 'User should be able to enter username'
$ws.Range("E72").Value = 'Click on "Search" button'
$ws.Range("E73").Value = 'Verify "Milyyy" record display'

# --- Step 2: New test case S9TC002 ("search an invalid system user") rows 74-77 ---
$ws.Range("C74").Value = 'S9TC002'
$ws.Range("D74").Value = 'User should be able to seach an invalid system user from admin page'
$ws.Range("E74").Value = 'Navigate to Admin Page'
$ws.Range("F74").Value = 'User should see Admin Page'
$ws.Range("E75").Value = 'Enter "Alex.X"in the "username" textField'
$ws.Range("F75").Value = 'User should be able to enter username'
$ws.Range("E76").Value = 'Click on "Search" button'
$ws.Range("F76").Value = 'User should be able to click "search" button'
$ws.Range("E77").Value = 'Verify "No Records Found" message display'
$ws.Range("F77").Value = 'User should see"No Records Found" messge'

# --- Step 3: Extend / add merges for the new rows ---
$ws.Range("A70:A77").Merge()
$ws.Range("B70:B77").Merge()
$ws.Range("C74:C77").Merge()
$ws.Range("D74:D77").Merge()

# --- Step 4: Styling ---
# D70:D73 switch to center/middle/wrap alignment (same as column D data cells)
$ws.Range("D70:D73").HorizontalAlignment = -4108
$ws.Range("D70:D73").VerticalAlignment = -4108
$ws.Range("D70:D73").WrapText = $true

# D74:D77 (new Test Case cell) center/middle/wrap alignment
$ws.Range("D74:D77").HorizontalAlignment = -4108
$ws.Range("D74:D77").VerticalAlignment = -4108
$ws.Range("D74:D77").WrapText = $true

# --- Step 5: Update selection / scroll position ---
$ws.Activate()
$ws.Range("E78").Select()

